$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45175 to 45183 for rows 2 through 13
$ws.Range("C2:C13").Value = 45183
